$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 175.8354644775391
$ws.Range("B1").Value = 204.3758087158203
$ws.Range("C1").Value = 53.14113235473633

$ws.Range("A3").Value = 173.7902374267578
$ws.Range("B3").Value = 191.04052734375
$ws.Range("C3").Value = 73.91013336181641

$ws.Range("A5").Value = 235.4834136962891
$ws.Range("B5").Value = 200.2127685546875
$ws.Range("C5").Value = 83.19190216064453

$ws.Range("A7").Value = 234.7833709716797
$ws.Range("B7").Value = 211.2177734375
$ws.Range("C7").Value = 66.07094573974609

$ws.Range("A9").Value = 238.6244354248047
$ws.Range("B9").Value = 182.9368133544922
$ws.Range("C9").Value = 77.99478912353516
